$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin stats (Price / Volume columns, and a few Coin/Link cells
# that were re-ordered). Each target cell is first forced to a text
# ("@") number format so that numeric-looking strings such as "218.70"
# or "24.01" are stored as text, matching the inline-string cells used
# in the source workbook instead of being auto-converted to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.327.92"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.717.03"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.87%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.70"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.44%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.01"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.77%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.08%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.959.77"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.715.30"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.59%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.43"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "28.303.21"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "249.77"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +5.83%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.43%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.58"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.62"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.56"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.53"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.44%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.67%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.479.68"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.02%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.17%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.970"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.41"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.19%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.99%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.59"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.60%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.65"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.76%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.863.99"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.53%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.27"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.22%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +6.66%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0114"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.25%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "90.03"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.79%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.06"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.81%  "
